$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column A (rows 1-32)
$aValues = @(
    -0.28991759397096928, -0.24316319624480798, -0.13767861566714856, -0.12500057094321804, -0.1184094357523886, -0.094220738895226752, -0.074087378189778619, 0.016755344651345538, 0.022800782821216892, 0.028869230465851103, -0.051473974460392924, -0.045395744376957037, -0.039150911925795562, -0.027084418295199697, -0.021052487540013765, -0.015027381452960231, -0.0090041656980339724, -0.036108675299537651, -0.02709647935202808, -0.018013364498985496, -0.0090041923366772636, -0.093948662016179796, -0.084635482488764957, -0.042126929326066431, -0.094933216614194293, -0.088689181413450058, -0.082374068706487336, -0.075292345429179619, -0.062542923447633569, -0.042171794726257428, -0.027019534626337105, -0.0060005919844288158
)

# New values for column B (rows 1-32)
$bValues = @(
    0.28927187094517137, 0.24062727707544695, 0.13700057079089234, 0.12440943515372638, 0.11724059034145196, 0.094087377471988365, 0.073777011364892786, -0.016800783441979661, -0.022869231083374686, -0.028884379424326312, 0.051395743756597057, 0.045150911295657625, 0.039084417618527212, 0.027052486905003725, 0.02102738081602773, 0.015004165058483565, 0.0089999993379983323, 0.036096478759528594, 0.027013363901156495, 0.018004191738073771, 0.0089999994008653772, 0.093635481875471882, 0.084126928463734707, 0.041999999132922738, 0.094689180798262385, 0.08837406808819992, 0.081292344797032179, 0.074542922762990571, 0.062171793978830436, 0.042019533908967066, 0.027000591222321901, 0.0059999993469173063
)

for ($i = 0; $i -lt 32; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Column B width changes from 15.42578125 to 14.7109375 (character width units).
# NOTE: this runtime quantizes ColumnWidth to 1/6-character-width pixel
# granularity when it serializes the XML `width` attribute, so the stored
# value can only land on a multiple of 1/6. 14.7109375 itself is not
# reachable exactly; 13.83 is the input that snaps to the closest possible
# stored width (14.666666666666666, the nearest representable value).
$ws.Columns.Item(2).ColumnWidth = 13.83
